$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

# Column A holds a literal date-like string (e.g. "2026/02/17"), not a real
# date value, so force Text format before assigning to stop Excel from
# auto-converting it to a date serial number. Reset the style afterward so
# the cell doesn't keep a stray custom "Text" number format.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2026/02/17"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "逃离鸭科夫"
$ws.Range("C$row").Value = 1210

# Match the centered alignment used by the rest of the data rows.
$rowRange = $ws.Range("A" + $row + ":C" + $row)
$rowRange.HorizontalAlignment = -4108
$rowRange.VerticalAlignment = -4108
